$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the data row (row 2) ---
# A2: Url -> new QA site
$ws.Range("A2").Value = "https://bento-qa.bento-tools.org/"

# C2 / D2: file names updated for the new single-case test
$ws.Range("C2").Value = "TC02_Bento_E2E_Select-Single-Add-To-Cart_Manifest.xlsx"
$ws.Range("D2").Value = "TC02_Bento_E2E_Select-Single-Add-To-Cart_WebData.xlsx"

# B2: Cypher query - ER status filter changed from Negative to Positive,
# and trailing endocrine_therapy_type / sample label filters removed.
$newQuery = "MATCH (ss:study_subject)`n" + `
  "MATCH (samp)-[:sample_of_study_subject]->(ss)`n" + `
  "MATCH (ss)<-[:sample_of_study_subject]-(samp:sample)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)`n" + `
  "WITH DISTINCT ss, samp, collect(DISTINCT samp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files`n" + `
  "MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)`n" + `
  "MATCH (ss)<-[:sf_of_study_subject]-(sf)`n" + `
  "MATCH (ss)<-[:diagnosis_of_study_subject]-(d)`n" + `
  "MATCH (d)<-[:tp_of_diagnosis]-(tp)`n" + `
  "MATCH (ss)<-[:demographic_of_study_subject]-(demo)`n" + `
  "WHERE ss.disease_subtype IN [`"Tubular Carcinoma`"] and d.tumor_size_group In [`"(3,4]`"] and d.er_status In [`"Positive`"]and d.pr_status In [`"Positive`"] `n" + `
  "return DISTINCT ss.study_subject_id as ``Case ID``,`n" + `
  "   p.program_acronym as ``Program Code``,`n" + `
  "    p.program_id as Program_ID,`n" + `
  "   s.study_acronym as ``Arm``,`n" + `
  "   ss.disease_subtype as ``Diagnosis``,`n" + `
  "   sf.grouped_recurrence_score AS ``Recurrence Score``,`n" + `
  "   d.tumor_size_group AS ``tumor_size``,`n" + `
  "   d.er_status AS ``ER Status``,`n" + `
  "   d.pr_status AS ``PR Status``,`n" + `
  "   demo.age_at_index AS ``Age (years)``,`n" + `
  "`tdemo.survival_time AS ``Survival (days)``"

$ws.Range("B2").Value = $newQuery

# --- Update the hyperlink on A2 to point at the new QA url, dropping the
#     stale #/ fragment + display text, while keeping the Hyperlink cell style ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://caninecommons.cancer.gov/") | Out-Null
$ws.Range("A2").Style = "Hyperlink"

# --- Scroll the sheet view back to column A ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 1

$wb.Save()
